# Update the NATMI LR-pair sheet (Ntf3-Ntrk1) with newly computed TPM values.
# Rows 5-7 of the old data are removed (the sending/target cluster combinations they
# described are consolidated into the remaining three rows), and rows 2-4 are updated
# with freshly calculated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete rows 5, 6 and 7 entirely, shrinking the sheet back down to
# a 4-row extent (header + 3 data rows).
$ws.Range("A5:T7").EntireRow.Delete()

# Row 2: ECs -> Ntf3/Ntrk1 -> FAPs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ntf3"
$ws.Range("C2").Value = "Ntrk1"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 9.736532333333333
$ws.Range("H2").Value = 29.209597
$ws.Range("I2").Value = 0.3545698647072128
$ws.Range("J2").Value = 0.3545698647072129
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.1211523333333333
$ws.Range("N2").Value = 0.363457
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 1.179603610758778
$ws.Range("R2").Value = 10.616432496829
$ws.Range("S2").Value = 0.3545698647072128
$ws.Range("T2").Value = 0.3545698647072129

# Row 3: FAPs -> Ntf3/Ntrk1 -> FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Ntf3"
$ws.Range("C3").Value = "Ntrk1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.17625966666667
$ws.Range("H3").Value = 30.528779
$ws.Range("I3").Value = 0.37058316962423
$ws.Range("J3").Value = 0.37058316962423
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.1211523333333333
$ws.Range("N3").Value = 0.363457
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 1.232877603222556
$ws.Range("R3").Value = 11.095898429003
$ws.Range("S3").Value = 0.37058316962423
$ws.Range("T3").Value = 0.37058316962423

# Row 4: MuSCs -> Ntf3/Ntrk1 -> FAPs
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Ntf3"
$ws.Range("C4").Value = "Ntrk1"
$ws.Range("D4").Value = "FAPs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 7.547331666666667
$ws.Range("H4").Value = 22.641995
$ws.Range("I4").Value = 0.2748469656685572
$ws.Range("J4").Value = 0.2748469656685572
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.1211523333333333
$ws.Range("N4").Value = 0.363457
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.9143768418572225
$ws.Range("R4").Value = 8.229391576715001
$ws.Range("S4").Value = 0.2748469656685572
$ws.Range("T4").Value = 0.2748469656685572
